$wb = $excel.ActiveWorkbook

function Set-Cell($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Clear-Cell($ws, $addr) {
    $ws.Range($addr).ClearContents()
}


# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 45
Set-Cell $ws "H45" 5000
Set-Cell $ws "J45" 5000
Set-Cell $ws "L45" 15000
Set-Cell $ws "N45" -15384
# Row 46
Set-Cell $ws "H46" 999
Set-Cell $ws "J46" 999
Set-Cell $ws "L46" 2997
Set-Cell $ws "N46" -3235
# Row 60
Set-Cell $ws "H60" 999
Set-Cell $ws "J60" 999
Set-Cell $ws "L60" 2997
Set-Cell $ws "N60" -3965

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 21
Set-Cell $ws "H21" 0
Set-Cell $ws "J21" 0
Set-Cell $ws "L21" 0
Clear-Cell $ws "N21"
# Row 63
Set-Cell $ws "H63" 3860.6667
Set-Cell $ws "I63" 3860.6667
Set-Cell $ws "J63" 0
Set-Cell $ws "K63" 3860.6667
Set-Cell $ws "L63" 0
Set-Cell $ws "M63" -3174.6667
Clear-Cell $ws "N63"
# Row 66
Set-Cell $ws "H66" 3860.6667
Set-Cell $ws "I66" 3860.6667
Set-Cell $ws "J66" 0
Set-Cell $ws "K66" 19303.3335
Set-Cell $ws "L66" 0
Set-Cell $ws "M66" -15871.3335
Clear-Cell $ws "N66"
# Row 76
Set-Cell $ws "H76" 39288
Set-Cell $ws "J76" 39288
Set-Cell $ws "L76" 39288
Set-Cell $ws "N76" -39964
# Row 79
Set-Cell $ws "H79" 39288
Set-Cell $ws "J79" 39288
Set-Cell $ws "L79" 39288
Set-Cell $ws "N79" -41628
# Row 110
Set-Cell $ws "H110" 0
Set-Cell $ws "I110" 0
Set-Cell $ws "J110" 0
Set-Cell $ws "K110" 0
Set-Cell $ws "L110" 0
Clear-Cell $ws "M110"
Clear-Cell $ws "N110"

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 19
Set-Cell $ws "H19" 0
Set-Cell $ws "J19" 0
Set-Cell $ws "L19" 0
Clear-Cell $ws "N19"
# Row 36
Set-Cell $ws "H36" 8405
Set-Cell $ws "I36" 4256.75
Set-Cell $ws "K36" 4256.75
Set-Cell $ws "M36" -3722.75
# Row 80
Set-Cell $ws "H80" 2688.5
Set-Cell $ws "I80" 2688.5
Set-Cell $ws "K80" 2688.5
Set-Cell $ws "M80" -1690.5
# Row 83
Set-Cell $ws "H83" 2688.5
Set-Cell $ws "I83" 2688.5
Set-Cell $ws "K83" 13442.5
Set-Cell $ws "M83" -8450.5

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
Set-Cell $ws "H7" 238.75
Set-Cell $ws "I7" 279.6
Set-Cell $ws "J7" 34.5
Set-Cell $ws "K7" 279.6
Set-Cell $ws "L7" 34.5
Set-Cell $ws "M7" -166.6
Set-Cell $ws "N7" -260.5
# Row 50
Set-Cell $ws "H50" 0
Set-Cell $ws "I50" 0
Set-Cell $ws "K50" 0
Clear-Cell $ws "M50"
# Row 93
Set-Cell $ws "H93" 45302.332
Set-Cell $ws "I93" 45302.332
Set-Cell $ws "K93" 45302.332
Set-Cell $ws "M93" -43430.332

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 2
Set-Cell $ws "H2" 740.13794
Set-Cell $ws "I2" 399.4
Set-Cell $ws "J2" 2869.75
Set-Cell $ws "K2" 2396.4
Set-Cell $ws "L2" 17218.5
Set-Cell $ws "M2" -2283.4
Set-Cell $ws "N2" -17444.5
# Row 23
Set-Cell $ws "H23" 0
Set-Cell $ws "J23" 0
Set-Cell $ws "L23" 0
Clear-Cell $ws "N23"
# Row 46
Set-Cell $ws "H46" 998
Set-Cell $ws "J46" 998
Set-Cell $ws "L46" 2994
Set-Cell $ws "N46" -3176
# Row 104
Set-Cell $ws "H104" 4591.75
Set-Cell $ws "J104" 4591.75
Set-Cell $ws "L104" 13775.25
Set-Cell $ws "N104" -19017.25
# Row 117
Set-Cell $ws "H117" 1914.9231
Set-Cell $ws "J117" 1657.8334
Set-Cell $ws "L117" 4973.5002
Set-Cell $ws "N117" -11857.5002

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 29
Set-Cell $ws "H29" 0
Set-Cell $ws "J29" 0
Set-Cell $ws "L29" 0
Clear-Cell $ws "N29"
# Row 102
Set-Cell $ws "H102" 5043.1665
Set-Cell $ws "I102" 4451.8
Set-Cell $ws "J102" 8000
Set-Cell $ws "K102" 4451.8
Set-Cell $ws "L102" 8000
Set-Cell $ws "M102" -2829.8
Set-Cell $ws "N102" -11244
# Row 107
Set-Cell $ws "H107" 293.75
Set-Cell $ws "J107" 0
Set-Cell $ws "L107" 0
Clear-Cell $ws "N107"
# Row 132
Set-Cell $ws "H132" 8399.5
Set-Cell $ws "I132" 8399.5
Set-Cell $ws "K132" 25198.5
Set-Cell $ws "M132" -22668.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 2
Set-Cell $ws "H2" 1001
Set-Cell $ws "I2" 1001
Set-Cell $ws "J2" 0
Set-Cell $ws "K2" 1001
Set-Cell $ws "L2" 0
Set-Cell $ws "M2" -889
Clear-Cell $ws "N2"
# Row 29
Set-Cell $ws "H29" 26600
Set-Cell $ws "I29" 0
Set-Cell $ws "J29" 26600
Set-Cell $ws "K29" 0
Set-Cell $ws "L29" 26600
Clear-Cell $ws "M29"
Set-Cell $ws "N29" -27190
# Row 43
Set-Cell $ws "H43" 25903.25
Set-Cell $ws "I43" 0
Set-Cell $ws "J43" 25903.25
Set-Cell $ws "K43" 0
Set-Cell $ws "L43" 25903.25
Clear-Cell $ws "M43"
Set-Cell $ws "N43" -26289.25
# Row 55
Set-Cell $ws "H55" 2269.7144
Set-Cell $ws "I55" 300
Set-Cell $ws "J55" 3057.6
Set-Cell $ws "K55" 300
Set-Cell $ws "L55" 3057.6
Set-Cell $ws "M55" -127
Set-Cell $ws "N55" -3403.6

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 34
Set-Cell $ws "H34" 10029
Set-Cell $ws "I34" 0
Set-Cell $ws "J34" 10029
Set-Cell $ws "K34" 0
Set-Cell $ws "L34" 10029
Clear-Cell $ws "M34"
Set-Cell $ws "N34" -10435
# Row 54
Set-Cell $ws "H54" 0
Set-Cell $ws "I54" 0
Set-Cell $ws "J54" 0
Set-Cell $ws "K54" 0
Set-Cell $ws "L54" 0
Clear-Cell $ws "M54"
Clear-Cell $ws "N54"
